$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.976.57"
$ws.Range("E2").Value = "  +2.98%  "

# Row 3
$ws.Range("D3").Value = "3.214.90"
$ws.Range("E3").Value = "  +1.86%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.17%  "

# Row 5
$ws.Range("D5").Value = "'604.46"
$ws.Range("E5").Value = "  +4.26%  "

# Row 6
$ws.Range("D6").Value = "'157.77"
$ws.Range("E6").Value = "  +5.52%  "

# Row 7
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.09%  "

# Row 8
$ws.Range("E8").Value = "  +6.34%  "

# Row 9
$ws.Range("D9").Value = "3.207.73"
$ws.Range("E9").Value = "  +1.64%  "

# Row 10
$ws.Range("D10").Value = "'0.162"
$ws.Range("E10").Value = "  +2.04%  "

# Row 11
$ws.Range("D11").Value = "'5.91"
$ws.Range("E11").Value = "  -3.72%  "

# Row 12
$ws.Range("E12").Value = "  +4.08%  "

# Row 13
$ws.Range("D13").Value = "'0.0000271"
$ws.Range("E13").Value = "  +2.45%  "

# Row 14
$ws.Range("D14").Value = "'39.53"
$ws.Range("E14").Value = "  +6.39%  "

# Row 15
$ws.Range("D15").Value = "3.732.82"
$ws.Range("E15").Value = "  +1.49%  "

# Row 16
$ws.Range("D16").Value = "66.772.79"
$ws.Range("E16").Value = "  +2.79%  "

# Row 17
$ws.Range("D17").Value = "'7.53"
$ws.Range("E17").Value = "  +5.54%  "

# Row 18
$ws.Range("D18").Value = "3.208.23"
$ws.Range("E18").Value = "  +1.37%  "

# Row 19
$ws.Range("D19").Value = "'527.33"
$ws.Range("E19").Value = "  +4.49%  "

# Row 20
$ws.Range("E20").Value = "  +0.81%  "

# Row 21
$ws.Range("D21").Value = "'15.58"
$ws.Range("E21").Value = "  +4.52%  "

# Row 22
$ws.Range("D22").Value = "'0.746"
$ws.Range("E22").Value = "  +4.40%  "

# Row 23
$ws.Range("D23").Value = "'8.27"
$ws.Range("E23").Value = "  +6.85%  "

# Row 24
$ws.Range("D24").Value = "'15.11"
$ws.Range("E24").Value = "  -0.67%  "

# Row 25
$ws.Range("D25").Value = "'85.91"
$ws.Range("E25").Value = "  +1.70%  "

# Row 26
$ws.Range("E26").Value = "  +0.04%  "

# Row 27
$ws.Range("D27").Value = "'9.31"
$ws.Range("E27").Value = "  +2.67%  "

# Row 28
$ws.Range("D28").Value = "'3.03"
$ws.Range("E28").Value = "  +3.54%  "

# Row 29
$ws.Range("E29").Value = "  +9.39%  "

# Row 30
$ws.Range("D30").Value = "'2.98"
$ws.Range("E30").Value = "  +6.95%  "

# Row 31
$ws.Range("D31").Value = "'7.06"
$ws.Range("E31").Value = "  +9.73%  "

# Row 32
$ws.Range("D32").Value = "'28.50"
$ws.Range("E32").Value = "  +3.34%  "

# Row 33
$ws.Range("E33").Value = "  +3.19%  "

# Row 34
$ws.Range("E34").Value = "  -0.01%  "

# Row 35
$ws.Range("D35").Value = "'6.61"
$ws.Range("E35").Value = "  +1.76%  "

# Row 36
$ws.Range("D36").Value = "'522.39"
$ws.Range("E36").Value = "  +9.25%  "

# Row 37
$ws.Range("D37").Value = "'54.95"
$ws.Range("E37").Value = "  -0.19%  "

# Row 38
$ws.Range("D38").Value = "'0.0912"
$ws.Range("E38").Value = "  +2.72%  "

# Row 39
$ws.Range("D39").Value = "'0.0428"
$ws.Range("E39").Value = "  +2.99%  "

# Row 40
$ws.Range("E40").Value = "  +8.68%  "

# Row 41
$ws.Range("D41").Value = "'8.95"
$ws.Range("E41").Value = "  +2.46%  "

# Row 42
$ws.Range("D42").Value = "'2.93"
$ws.Range("E42").Value = "  -0.18%  "

# Row 43
$ws.Range("D43").Value = "0.0₃0687"
$ws.Range("E43").Value = "  +16.12%  "

# Row 44
$ws.Range("D44").Value = "'0.303"
$ws.Range("E44").Value = "  +7.25%  "

# Row 45
$ws.Range("E45").Value = "  +1.78%  "

# Row 46
$ws.Range("D46").Value = "2.905.52"
$ws.Range("E46").Value = "  -2.86%  "

# Row 47
$ws.Range("D47").Value = "'28.76"
$ws.Range("E47").Value = "  +1.68%  "

# Row 48
$ws.Range("D48").Value = "'2.76"
$ws.Range("E48").Value = "  +11.38%  "

# Row 49
$ws.Range("E49").Value = "  +4.15%  "

# Row 50
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").Value = "'2.36"
$ws.Range("E50").Value = "  +5.10%  "

# Row 51
$ws.Range("B51").Value = "USDe"
$ws.Range("C51").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D51").Value = "'0.999"
$ws.Range("E51").Value = "  -0.04%  "
